$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A13").Value = 44455
$ws.Range("B13").Value = "model_floodwater_unet_pc_augm_diceloss_without_square"
$ws.Range("C13").Value = 0.698

$ws.Range("A14").Value = 44456
$ws.Range("B14").Value = "model_floodwater_unet_pc_augm_diceloss_2"
$ws.Range("E14").Value = "['hbe', 'jja']"

$ws.Range("E15").Select()
